$wb = $excel.ActiveWorkbook

# --- Sheet "Đơn phụ phẫu 1": collapse columns J:AA and repurpose G/H/I ---
$ws2 = $wb.Worksheets.Item(2)

# Remove all the now-unused trailing columns (J..AA) for rows 1-3
$ws2.Range("J1:AA3").ClearContents()

# Row 1 headers shift meaning: G/H/I get new labels
$ws2.Range("G1").Value = "Tên dịch vụ"
$ws2.Range("H1").Value = "Phụ phẫu 1"
$ws2.Range("I1").Value = "Công phụ phẫu 1"

# Row 2 data values
$ws2.Range("G2").Value = "cắt sẹo "
$ws2.Range("H2").Value = "Trần Khánh Hiệp"
$ws2.Range("I2").ClearContents()

# Row 3 totals
$ws2.Range("I3").Value = 0

# --- Sheet "Lương": rename label and update computed totals ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A1").Value = "Danh mục lương"
$ws3.Range("B30").Value = 1550000
$ws3.Range("B31").Value = 1550000
